# Updated symbol list on Thu Jan 26 07:36:06 UTC 2023 with GitHub Actions
# Refresh crypto price/volume figures in the "cryptos" sheet.
# Values are stored as literal text in the sheet (e.g. "307.59", "1.40%"),
# so each new value is written with a leading apostrophe to force a text
# entry and avoid Excel auto-converting it to a number/percentage.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'307.59"
$ws.Range("E2").Value = "'1.40%"
$ws.Range("E3").Value = "'1.44%"
$ws.Range("D4").Value = "'5.110"
$ws.Range("E4").Value = "'0.92%"
$ws.Range("D5").Value = "'0.08090"
$ws.Range("E5").Value = "'1.07%"
$ws.Range("D6").Value = "'1.953"
$ws.Range("E6").Value = "'1.60%"
$ws.Range("D7").Value = "'4.198"
$ws.Range("E7").Value = "'3.70%"
$ws.Range("D8").Value = "'7.758"
$ws.Range("E8").Value = "'0.14%"
$ws.Range("D9").Value = "'0.9287"
$ws.Range("E9").Value = "'0.82%"
$ws.Range("D10").Value = "'0.1380"
$ws.Range("E10").Value = "'13.26%"
$ws.Range("D11").Value = "'0.1919"
$ws.Range("E11").Value = "'3.80%"
$ws.Range("D12").Value = "'0.09219"
$ws.Range("E12").Value = "'-2.29%"
$ws.Range("D13").Value = "'0.03454"
$ws.Range("E13").Value = "'-3.20%"
$ws.Range("D14").Value = "'0.09845"
$ws.Range("E14").Value = "'-0.10%"
$ws.Range("D15").Value = "'0.001419"
$ws.Range("E15").Value = "'2.22%"
$ws.Range("D16").Value = "'0.005807"
$ws.Range("E16").Value = "'0.31%"
$ws.Range("E17").Value = "'3.71%"
$ws.Range("D18").Value = "'3.006"
$ws.Range("E18").Value = "'2.14%"
$ws.Range("D19").Value = "'0.3447"
$ws.Range("E19").Value = "'1.02%"
$ws.Range("D20").Value = "'0.1315"
$ws.Range("E20").Value = "'2.57%"
$ws.Range("D21").Value = "'4.910"
$ws.Range("E21").Value = "'-2.61%"
$ws.Range("D22").Value = "'0.2445"
$ws.Range("E22").Value = "'-0.79%"
$ws.Range("D23").Value = "'0.04453"
$ws.Range("E23").Value = "'-1.24%"
$ws.Range("D24").Value = "'0.001220"
$ws.Range("E24").Value = "'0.34%"
$ws.Range("D25").Value = "'0.004830"
$ws.Range("E25").Value = "'-0.38%"
$ws.Range("D26").Value = "'0.0001243"
$ws.Range("E26").Value = "'-0.57%"
$ws.Range("D39").Value = "'0.02018"
$ws.Range("E39").Value = "'4.43%"
$ws.Range("D40").Value = "'0.04930"
$ws.Range("E40").Value = "'3.74%"
$ws.Range("D41").Value = "'0.007680"
$ws.Range("E41").Value = "'2.15%"
$ws.Range("D42").Value = "'0.01011"
$ws.Range("E42").Value = "'5.84%"
$ws.Range("D43").Value = "'0.1376"
$ws.Range("E43").Value = "'3.36%"
$ws.Range("D44").Value = "'0.002105"
$ws.Range("E44").Value = "'-0.25%"
$ws.Range("D45").Value = "'0.01160"
$ws.Range("E45").Value = "'3.83%"
$ws.Range("D46").Value = "'0.00006455"
$ws.Range("E46").Value = "'2.67%"
$ws.Range("E47").Value = "'0.23%"
$ws.Range("D48").Value = "'63.57"
$ws.Range("E48").Value = "'-1.41%"
$ws.Range("E49").Value = "'-8.54%"
$ws.Range("E50").Value = "'0.23%"
$ws.Range("E51").Value = "'0.23%"
